# Commit: "Add New TC for Sim ATM"
# Rename the worksheet to match the new test case id and move the
# active selection to where the new test-case row is being added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet name: Simatm14_1_1_3_BH -> Simatm14_1_1_1_BH
$ws.Name = "Simatm14_1_1_1_BH"

# Move / record the active selection at D20 (where the new TC row starts)
$ws.Range("D20").Select()
